$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated PSSM score matrix (B2:K21) with supplemental-figure recalculated values
$data = @(
    @(-20.37389107519477, 2.451346657547671, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477),
    @(-20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, 2.482560480274871, -20.37389107519477, -20.37389107519477),
    @(-20.37389107519477, 2.182475815275993, 2.837553014622893, -20.37389107519477, 2.574633898330255, -20.37389107519477, 1.914075162532818, -20.37389107519477, 2.385082141610652, -20.37389107519477),
    @(-20.37389107519477, 0.9757060239461914, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477),
    @(-20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477),
    @(2.956030737821884, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477),
    @(-20.37389107519477, -20.37389107519477, -20.37389107519477, 2.872669923164519, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477),
    @(3.613546520808597, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477),
    @(-20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, 1.537476406966603, -20.37389107519477, 2.435535781179276),
    @(-20.37389107519477, -20.37389107519477, -20.37389107519477, 2.022136921138051, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, 1.594807002712864),
    @(-20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477),
    @(-20.37389107519477, -20.37389107519477, -20.37389107519477, 1.698011495781001, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, 2.158432743211956, 1.514418520679353),
    @(-20.37389107519477, -20.37389107519477, 1.696852591715197, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, 1.894640930975546),
    @(-20.37389107519477, -20.37389107519477, -0.1730008492572386, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477),
    @(-20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, 2.274049707254667, -20.37389107519477),
    @(-20.37389107519477, 0.8210611101204024, 0.07421817066045892, -20.37389107519477, -20.37389107519477, -20.37389107519477, 0.7460537300147408, 1.040856604208673, 1.326730245195369, -20.37389107519477),
    @(-20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, -20.37389107519477, 0.6467016743409724, 1.166061433845963, 1.568923001207022, -20.37389107519477),
    @(-20.37389107519477, -20.37389107519477, 1.73369414036787, -20.37389107519477, -20.37389107519477, -20.37389107519477, 1.750408545117687, 1.987625742019606, -20.37389107519477, -20.37389107519477),
    @(-20.37389107519477, 1.604043173053537, 2.119165183828499, -20.37389107519477, 3.811757568652342, -20.37389107519477, 2.130795308548909, 1.696426771660264, -20.37389107519477, 2.3203879481971),
    @(-20.37389107519477, 1.686235189027892, -20.37389107519477, 2.424763647456333, -20.37389107519477, 4.321927086224679, 2.390766130948589, -20.37389107519477, -20.37389107519477, -20.37389107519477)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}
